# Updated cryptos list - price (D) and 1h volume-change (E) refresh,
# plus a swap of two row pairs (Uniswap/BitcoinCash, WhiteBITCoin/ImmutableX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.270.35'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '2.642.56'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.14'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.24'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('E9').Value = '  +5.02%  '
$ws.Range('D10').Value = '2.641.97'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.360'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.02%  '
$ws.Range('E14').Value = '  +3.63%  '
$ws.Range('D15').Value = '3.111.90'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '72.156.43'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').Value = '2.648.61'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.99'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.43%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.94'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.42'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.09'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +11.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.39'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.41'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.13'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.07%  '
$ws.Range('D28').Value = '2.782.72'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').Value = '0.0₃0960'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '525.74'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.34'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.41'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('E38').Value = '  -5.86%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.08'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.40'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.86'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.65'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.10'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.37'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '151.46'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.73'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('E51').Value = '  -3.85%  '
